# Apply the "Assets" sheet changes: add two new asset rows
# (HashURL, ACME_WorkItemsUpdateURL) used by the new
# System1_Extract_ClientInformation / System1_Update_WorkItem workflows,
# and update the sheet's active selection/scroll position.

$wb = $excel.ActiveWorkbook

# The sheet that changed is "Assets" (3rd tab in the workbook).
$ws = $wb.Worksheets.Item(3)

# Row 6: new "HashURL" asset entry (Name / Asset / Env / Description columns)
$ws.Range("A6").Value = "HashURL"
$ws.Range("B6").Value = "HashURL"
$ws.Range("C6").Value = "Dev"
$ws.Range("D6").Value = "HashURL"

# Row 7: new "ACME_WorkItemsUpdateURL" asset entry
$ws.Range("A7").Value = "ACME_WorkItemsUpdateURL"
$ws.Range("B7").Value = "ACME_WorkItemsUpdateURL"
$ws.Range("C7").Value = "Dev"
$ws.Range("D7").Value = "ACME_WorkItemsUpdateURL"

# Reflect the updated view state: sheet scrolled right one column and the
# selected cell moved to C8.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("C8").Select()
